# Modificaciones por error en reporte de listado de control y en migracion de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark "Error en listado de control - pagos de la fecha" (row 63) as done (100%)
$ws.Range("C63").Value = 1
$ws.Range("C63").NumberFormat = "0%"

# Mark "Migracion de datos" (row 64) as done (100%)
$ws.Range("C64").Value = 1
$ws.Range("C64").NumberFormat = "0%"

# Mark "Reporte para contador" (row 67) as "en proceso"
$ws.Range("C67").Value = "en proceso"

# Update view state to match author's scroll/selection position
$excel.ActiveWindow.ScrollRow = 52
$null = $ws.Range("C68").Select()
